# MYR refinements (#334)
#
# 1) On the "Details" sheet, drop the "Supplier Name" and "Model Year"
#    columns and reorder what's left to: Classification, Compliant,
#    ZEV Class Ordering.
# 2) Remove the now-unused "Divisor" worksheet.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$details = $wb.Worksheets.Item("Details")

# Swap "ZEV Class Ordering" (D1) and "Compliant" (E1) first so that once
# the leading two columns are removed, the remaining three land in the
# target order: Classification | Compliant | ZEV Class Ordering.
$tmp = $details.Range("D1").Value2
$details.Range("D1").Value2 = $details.Range("E1").Value2
$details.Range("E1").Value2 = $tmp

# Drop "Supplier Name" (A1) and "Model Year" (B1), shifting the rest left.
$details.Range("A1:B1").Delete(-4159)

# Match the refreshed column widths (~24.5 / 17.5 / 23.16 / 17 characters).
$details.Columns.Item(1).ColumnWidth = 23.666666666666668
$details.Columns.Item(2).ColumnWidth = 16.666666666666668
$details.Columns.Item(3).ColumnWidth = 22.33
$details.Columns.Item(4).ColumnWidth = 16.166666666666668

# Remove the "Divisor" worksheet entirely.
$wb.Worksheets.Item("Divisor").Delete()

# Keep "Details" the active/selected tab, as it was originally.
$details.Activate()
$details.Select()
